$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is reported for this market/product, so it is
# inserted as the new row 3 (most recent date) and every existing
# data row shifts down by one (old row 3 -> new row 4, ..., old row 23 -> new row 24).
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44496
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 84
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = 1800
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1800
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
